$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 5 (G5 / H5): drop the now-unused TakeScreenshot(VT200_0605...) steps
#     and swap the Screenshot-based validations for isIconDisplayed checks.
$ws.Range("H5").Value = "validate1`n{`nvalidate_PageTitle=Compliance JS specs`n};`nvalidate2`n{`nvalidate_PageTitle=Native Toolbar JS Test`n};`nvalidate3`n{`nvalidate_Text_Exists=VT200-0603`n};`nvalidate4`n{`nvalidate_isIconDisplayed=toobarview_xpath,true`n};`nvalidate5`n{`nvalidate_Text_Exists=VT200-0605`n};`nvalidate6`n{`nvalidate_isIconDisplayed=toobarview_xpath,false`n};`n"

$ws.Range("G5").Value = "wait(5);`nvalidate1;`nlink_Click(toolbar_test_link);`nvalidate2;`nSelectTestToRun(VT200_0603_string);`nClickRunTest(runtest_top_xpath);`nvalidate3;`nClickRunTest(runtest_bottom_xpath);`nwait(5);`nvalidate4;`nSelectTestToRun(VT200_0605_string);`nClickRunTest(runtest_top_xpath);`nvalidate5;`nClickRunTest(runtest_bottom_xpath);`nwait(5);`nvalidate6;`n"

# --- Row 10 / D10: pick up the same (bottom-bordered) formatting the rest of
#     column D already uses, matching the rest of the table's last row.
$ws.Range("D2").Copy() | Out-Null
$ws.Range("D10").PasteSpecial(-4122) | Out-Null

# --- Row 10 height + view selection/scroll position.
$ws.Rows.Item(10).RowHeight = 409.6

$excel.CutCopyMode = $false

$ws.Range("G2").Select()
$ws.Application.ActiveWindow.ScrollRow = 1
